$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 114, pushing existing rows 114:134 down to 115:135
$ws.Rows.Item(114).Insert()

# Populate the new row 114 with the same style/formatting pattern as its
# neighbours and the values specified by the diff.
$ws.Cells.Item(114, 1).Value = 5
$ws.Cells.Item(114, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(114, 3).Value = "Maule"
$ws.Cells.Item(114, 4).Value = 45209
$ws.Cells.Item(114, 4).Style = $ws.Cells.Item(115, 4).Style
$ws.Cells.Item(114, 4).NumberFormat = $ws.Cells.Item(115, 4).NumberFormat
$ws.Cells.Item(114, 5).Value = 7
$ws.Cells.Item(114, 6).Value = 100112040
$ws.Cells.Item(114, 7).Value = "Cilantro"
$ws.Cells.Item(114, 8).Value = "Sin especificar"
$ws.Cells.Item(114, 9).Value = "Primera"
$ws.Cells.Item(114, 10).Value = 150
$ws.Cells.Item(114, 11).Value = 9000
$ws.Cells.Item(114, 12).Value = 9000
$ws.Cells.Item(114, 13).Value = 9000
$ws.Cells.Item(114, 14).Value = "`$/caja 36 atados"
$ws.Cells.Item(114, 15).Value = "Región Metropolitana"
$ws.Cells.Item(114, 16).Value = 250
$ws.Cells.Item(114, 17).Value = 36
$ws.Cells.Item(114, 18).Value = "Hortaliza"
